$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Range("D2").Value = "28.492.29"
$ws.Range("E2").Value = "  +2.40%  "
$ws.Range("D3").Value = "1.829.56"
$ws.Range("E3").Value = "  +2.15%  "
Set-TextValue "D4" "1.001"
$ws.Range("E4").Value = "  +0.01%  "
Set-TextValue "D5" "316.08"
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("E6").Value = "  +0.05%  "
Set-TextValue "D7" "0.5074"
$ws.Range("E7").Value = "  -4.77%  "
Set-TextValue "D8" "0.3914"
$ws.Range("E8").Value = "  +2.10%  "
Set-TextValue "D9" "0.07710"
$ws.Range("E9").Value = "  +3.93%  "
$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue "D10" "1.116"
$ws.Range("E10").Value = "  +3.06%  "
$ws.Range("B11").Value = "OKB"
$ws.Range("C11").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D11" "41.93"
$ws.Range("E11").Value = "  +1.52%  "
Set-TextValue "D12" "21.06"
$ws.Range("E12").Value = "  +3.94%  "
Set-TextValue "D13" "6.282"
$ws.Range("E13").Value = "  +1.72%  "
Set-TextValue "D14" "7.577"
$ws.Range("E14").Value = "  +1.99%  "
$ws.Range("E15").Value = "  +0.05%  "
$ws.Range("D16").Value = "1.826.57"
$ws.Range("E16").Value = "  +1.85%  "
Set-TextValue "D17" "93.55"
$ws.Range("E17").Value = "  +6.32%  "
Set-TextValue "D18" "0.00001085"
$ws.Range("E18").Value = "  +2.71%  "
Set-TextValue "D19" "0.06639"
$ws.Range("E19").Value = "  +1.86%  "
Set-TextValue "D20" "17.72"
$ws.Range("E20").Value = "  +3.04%  "
$ws.Range("E21").Value = "  +0.09%  "
Set-TextValue "D22" "6.159"
$ws.Range("E22").Value = "  +3.62%  "
$ws.Range("D23").Value = "28.517.50"
$ws.Range("E23").Value = "  +2.31%  "
Set-TextValue "D24" "11.12"
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("E25").Value = "  +7.80%  "
Set-TextValue "D26" "156.97"
$ws.Range("E26").Value = "  -0.01%  "
Set-TextValue "D27" "20.65"
$ws.Range("E27").Value = "  +2.82%  "
$ws.Range("D28").Value = "2.035.88"
$ws.Range("E28").Value = "  +1.90%  "
Set-TextValue "D29" "2.412"
$ws.Range("E29").Value = "  +4.71%  "
Set-TextValue "D30" "125.41"
$ws.Range("E30").Value = "  +3.49%  "
Set-TextValue "D31" "1.137"
$ws.Range("E31").Value = "  +3.71%  "
$ws.Range("E32").Value = "  -0.17%  "
Set-TextValue "D33" "5.674"
$ws.Range("E33").Value = "  +3.37%  "
$ws.Range("E34").Value = "  +0.27%  "
Set-TextValue "D35" "0.07068"
$ws.Range("E35").Value = "  +2.19%  "
Set-TextValue "D36" "0.2225"
$ws.Range("E36").Value = "  +1.32%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D37" "0.02329"
$ws.Range("E37").Value = "  +2.94%  "
$ws.Range("B38").Value = "FraxShare"
$ws.Range("C38").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D38" "8.893"
$ws.Range("E38").Value = "  +6.51%  "
Set-TextValue "D39" "5.158"
$ws.Range("E39").Value = "  +2.48%  "
Set-TextValue "D40" "0.6257"
$ws.Range("E40").Value = "  +2.81%  "
Set-TextValue "D41" "11.20"
$ws.Range("E41").Value = "  -0.62%  "
Set-TextValue "D42" "1.184"
$ws.Range("E42").Value = "  +0.87%  "
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("E44").Value = "  -1.00%  "
Set-TextValue "D45" "13.49"
$ws.Range("E45").Value = "  +2.05%  "
Set-TextValue "D46" "0.5906"
$ws.Range("E46").Value = "  +3.87%  "
Set-TextValue "D47" "3.717"
$ws.Range("E47").Value = "  +0.99%  "
Set-TextValue "D48" "124.90"
$ws.Range("E48").Value = "  +0.63%  "
Set-TextValue "D49" "1.978"
$ws.Range("E49").Value = "  +3.66%  "
Set-TextValue "D50" "1.197"
$ws.Range("E50").Value = "  +2.57%  "
Set-TextValue "D51" "0.06925"
$ws.Range("E51").Value = "  +1.94%  "
